# Regenerate orders with updated distance/size codes.
# The experiment's distance and size conditions were renamed:
#   D80 -> D86, D64 -> D69, D51 -> D55, S30 -> S31
# These tokens appear embedded throughout the Condition / Filename_Left /
# Filename_Right / Distance / Size columns (e.g. "Face14_D80_S20",
# "Fixation_D80_l.png", "D80", "S30"), so apply a global find & replace
# across the whole sheet for each of the four token substitutions.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Replace("D80", "D86")
$ws.Cells.Replace("D64", "D69")
$ws.Cells.Replace("D51", "D55")
$ws.Cells.Replace("S30", "S31")
